$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 168.04
$ws.Range("I15").Value = 168.04
$ws.Range("K15").Value = 504.12
$ws.Range("M15").Value = -335.12

$ws.Range("H62").Value = 2500.4614
$ws.Range("I62").Value = 2713.25
$ws.Range("K62").Value = 2713.25
$ws.Range("M62").Value = -2089.25

$ws.Range("H65").Value = 2500.4614
$ws.Range("I65").Value = 2713.25
$ws.Range("K65").Value = 13566.25
$ws.Range("M65").Value = -10446.25

$ws.Range("H112").Value = 58824920
$ws.Range("J112").Value = 71430136
$ws.Range("L112").Value = 214290408
$ws.Range("N112").Value = -214292624

$ws.Range("H129").Value = 914.6
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 914.6
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 2743.8
$ws.Range("M129").ClearContents()
$ws.Range("N129").Value = -12743.8

$ws.Range("H138").Value = 2783.6572
$ws.Range("I138").Value = 2871.4666
$ws.Range("J138").Value = 2717.8
$ws.Range("K138").Value = 8614.399800000001
$ws.Range("L138").Value = 8153.400000000001
$ws.Range("M138").Value = -3474.399800000001
$ws.Range("N138").Value = -18433.4

$ws.Range("H140").Value = 49713.332
$ws.Range("J140").Value = 49713.332
$ws.Range("L140").Value = 49713.332
$ws.Range("N140").Value = -60073.332

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6630
$ws.Range("I32").Value = 6072.9165
$ws.Range("J32").Value = 20000
$ws.Range("K32").Value = 6072.9165
$ws.Range("L32").Value = 20000
$ws.Range("M32").Value = -5785.9165
$ws.Range("N32").Value = -20574

$ws.Range("H74").Value = 8242.6
$ws.Range("I74").Value = 1327.1538
$ws.Range("J74").Value = 21085.572
$ws.Range("K74").Value = 1327.1538
$ws.Range("L74").Value = 21085.572
$ws.Range("M74").Value = -453.1538
$ws.Range("N74").Value = -22833.572

$ws.Range("H77").Value = 8242.6
$ws.Range("I77").Value = 1327.1538
$ws.Range("J77").Value = 21085.572
$ws.Range("K77").Value = 6635.769
$ws.Range("L77").Value = 105427.86
$ws.Range("M77").Value = -2267.769
$ws.Range("N77").Value = -114163.86

$ws.Range("H132").Value = 14974169
$ws.Range("I132").Value = 25596700
$ws.Range("J132").Value = 102626
$ws.Range("K132").Value = 76790100
$ws.Range("L132").Value = 307878
$ws.Range("M132").Value = -76787570
$ws.Range("N132").Value = -312938

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 48914
$ws.Range("J59").Value = 48914
$ws.Range("L59").Value = 48914
$ws.Range("N59").Value = -50608

$ws.Range("H94").Value = 480.42856
$ws.Range("I94").Value = 427.16666
$ws.Range("J94").Value = 800
$ws.Range("K94").Value = 427.16666
$ws.Range("L94").Value = 800
$ws.Range("M94").Value = 23.83334000000002
$ws.Range("N94").Value = -1702

$ws.Range("H107").Value = 6567.304
$ws.Range("I107").Value = 5859.25
$ws.Range("J107").Value = 8185.7144
$ws.Range("K107").Value = 5859.25
$ws.Range("L107").Value = 8185.7144
$ws.Range("M107").Value = -3939.25
$ws.Range("N107").Value = -12025.7144

$ws.Range("H134").Value = 116581.3
$ws.Range("I134").Value = 116581.3
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 349743.9
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -347208.9
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 8000
$ws.Range("J17").Value = 8000
$ws.Range("L17").Value = 8000
$ws.Range("N17").Value = -8348

$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").ClearContents()
$ws.Range("N25").ClearContents()

$ws.Range("H31").Value = 9808580
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 9808580
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 9808580
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -9809170

$ws.Range("H34").Value = 9808580
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 9808580
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 9808580
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -9808984

$ws.Range("H41").Value = 15600
$ws.Range("J41").Value = 19933.334
$ws.Range("L41").Value = 19933.334
$ws.Range("N41").Value = -20789.334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 142.42105
$ws.Range("J23").Value = 109.35714
$ws.Range("L23").Value = 328.07142
$ws.Range("N23").Value = -798.07142

$ws.Range("H131").Value = 2415.7014
$ws.Range("J131").Value = 1586.2812
$ws.Range("L131").Value = 4758.8436
$ws.Range("N131").Value = -14838.8436

$ws.Range("H133").Value = 4975
$ws.Range("I133").Value = 4971.4287
$ws.Range("J133").Value = 5000
$ws.Range("K133").Value = 14914.2861
$ws.Range("L133").Value = 15000
$ws.Range("M133").Value = -9854.286100000001
$ws.Range("N133").Value = -25120

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4554.919
$ws.Range("I70").Value = 4386.8887
$ws.Range("J70").Value = 4714.1055
$ws.Range("K70").Value = 4386.8887
$ws.Range("L70").Value = 4714.1055
$ws.Range("M70").Value = -4116.8887
$ws.Range("N70").Value = -5254.1055

$ws.Range("H73").Value = 4554.919
$ws.Range("I73").Value = 4386.8887
$ws.Range("J73").Value = 4714.1055
$ws.Range("K73").Value = 4386.8887
$ws.Range("L73").Value = 4714.1055
$ws.Range("M73").Value = -3450.8887
$ws.Range("N73").Value = -6586.1055

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 5750
$ws.Range("J42").Value = 8000
$ws.Range("L42").Value = 8000
$ws.Range("N42").Value = -9126

$ws.Range("H49").Value = 5750
$ws.Range("J49").Value = 8000
$ws.Range("L49").Value = 8000
$ws.Range("N49").Value = -8294

$ws.Range("H132").Value = 7042.96
$ws.Range("I132").Value = 7898.8
$ws.Range("J132").Value = 3619.6
$ws.Range("K132").Value = 23696.4
$ws.Range("L132").Value = 10858.8
$ws.Range("M132").Value = -21166.4
$ws.Range("N132").Value = -15918.8

$ws.Range("H136").Value = 1679.75
$ws.Range("I136").Value = 1011.5263
$ws.Range("J136").Value = 4219
$ws.Range("K136").Value = 3034.5789
$ws.Range("L136").Value = 12657
$ws.Range("M136").Value = -484.5789
$ws.Range("N136").Value = -17757

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4264.091
$ws.Range("I62").Value = 4217
$ws.Range("J62").Value = 4320.6
$ws.Range("K62").Value = 4217
$ws.Range("L62").Value = 4320.6
$ws.Range("M62").Value = -3593
$ws.Range("N62").Value = -5568.6

$ws.Range("H65").Value = 4264.091
$ws.Range("I65").Value = 4217
$ws.Range("J65").Value = 4320.6
$ws.Range("K65").Value = 21085
$ws.Range("L65").Value = 21603
$ws.Range("M65").Value = -17965
$ws.Range("N65").Value = -27843

$ws.Range("H136").Value = 3670692.5
$ws.Range("I136").Value = 12221.546
$ws.Range("J136").Value = 8405184
$ws.Range("K136").Value = 36664.638
$ws.Range("L136").Value = 25215552
$ws.Range("M136").Value = -34114.638
$ws.Range("N136").Value = -25220652
